$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full updated data table (rows 2-13), columns A:T
# Columns: A=Sending cluster, B=Ligand symbol, C=Receptor symbol, D=Target cluster,
# E=Ligand-expressing cells, F=Ligand detection rate, G=Ligand average expression value,
# H=Ligand total expression value, I=Ligand derived specificity of average expression value,
# J=Ligand derived specificity of total expression value, K=Receptor-expressing cells,
# L=Receptor detection rate, M=Receptor average expression value, N=Receptor total expression value,
# O=Receptor derived specificity of average expression value, P=Receptor derived specificity of total expression value,
# Q=Edge average expression weight, R=Edge total expression weight,
# S=Edge average expression derived specificity, T=Edge total expression derived specificity

$data = @(
    @("ECs","Ntrk3","Ptprs","ECs",1,0.3333333333333333,0.004481,0.013443,0.01111165297720057,0.01111165297720057,3,1,5.273684,15.821052,0.0510821201937383,0.0510821201937383,0.023631378004,0.212682402036,0.0005676067929324695,0.0005676067929324696),
    @("ECs","Ntrk3","Ptprs","FAPs",1,0.3333333333333333,0.004481,0.013443,0.01111165297720057,0.01111165297720057,3,1,57.80210733333333,173.406322,0.5598845502029881,0.5598845502029881,0.2590112429606666,2.331101186646,0.006221242829151635,0.006221242829151635),
    @("ECs","Ntrk3","Ptprs","MuSCs",1,0.3333333333333333,0.004481,0.013443,0.01111165297720057,0.01111165297720057,3,1,32.95839133333334,98.87517400000002,0.3192425840231603,0.3192425840231604,0.1476865515646667,1.329178964082,0.003547312809210152,0.003547312809210154),
    @("ECs","Ntrk3","Ptprs","Resolving-Mac",1,0.3333333333333333,0.004481,0.013443,0.01111165297720057,0.01111165297720057,3,1,7.205150000000001,21.61545,0.06979074558011317,0.06979074558011318,0.03228627715,0.29057649435,0.000775490545906312,0.0007754905459063122),
    @("FAPs","Ntrk3","Ptprs","ECs",3,1,0.1511493333333333,0.453448,0.3748089577628241,0.3748089577628241,3,1,5.273684,15.821052,0.0510821201937383,0.0510821201937383,0.7971138208106667,7.174024387295999,0.01914603623013036,0.01914603623013036),
    @("FAPs","Ntrk3","Ptprs","FAPs",3,1,0.1511493333333333,0.453448,0.3748089577628241,0.3748089577628241,3,1,57.80210733333333,173.406322,0.5598845502029881,0.5598845502029881,8.736749988695111,78.63074989825598,0.2098497447290895,0.2098497447290895),
    @("FAPs","Ntrk3","Ptprs","MuSCs",3,1,0.1511493333333333,0.453448,0.3748089577628241,0.3748089577628241,3,1,32.95839133333334,98.87517400000002,0.3192425840231603,0.3192425840231604,4.981638877772445,44.834749899952,0.1196549801912315,0.1196549801912315),
    @("FAPs","Ntrk3","Ptprs","Resolving-Mac",3,1,0.1511493333333333,0.453448,0.3748089577628241,0.3748089577628241,3,1,7.205150000000001,21.61545,0.06979074558011317,0.06979074558011318,1.089053619066667,9.801482571600001,0.02615819661237264,0.02615819661237264),
    @("MuSCs","Ntrk3","Ptprs","ECs",3,1,0.24764,0.74292,0.6140793892599753,0.6140793892599754,3,1,5.273684,15.821052,0.0510821201937383,0.0510821201937383,1.30597510576,11.75377595184,0.03136847717067547,0.03136847717067547),
    @("MuSCs","Ntrk3","Ptprs","FAPs",3,1,0.24764,0.74292,0.6140793892599753,0.6140793892599754,3,1,57.80210733333333,173.406322,0.5598845502029881,0.5598845502029881,14.31411386002667,128.82702474024,0.3438135626447469,0.343813562644747),
    @("MuSCs","Ntrk3","Ptprs","MuSCs",3,1,0.24764,0.74292,0.6140793892599753,0.6140793892599754,3,1,32.95839133333334,98.87517400000002,0.3192425840231603,0.3192425840231604,8.161816029786667,73.45634426808002,0.1960402910227186,0.1960402910227187),
    @("MuSCs","Ntrk3","Ptprs","Resolving-Mac",3,1,0.24764,0.74292,0.6140793892599753,0.6140793892599754,3,1,7.205150000000001,21.61545,0.06979074558011317,0.06979074558011318,1.784283346,16.058550114,0.04285705842183422,0.04285705842183423)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $startRow + $i
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Count; $j++) {
        $col = $j + 1
        $ws.Cells.Item($rowNum, $col).Value = $rowValues[$j]
    }
}

$wb.Save()
